$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 34.32219813181192
$ws.Range("G2").Value = 33.93347625483653
$ws.Range("H2").Value = 34.70490756899319
$ws.Range("I2").Value = 3.833097380741688
$ws.Range("J2").Value = 3.805469405110184
$ws.Range("K2").Value = 3.86051984770603
$ws.Range("L2").Value = 0.2688359171456197
$ws.Range("M2").Value = 0.2667845704221329
$ws.Range("N2").Value = 0.2708811391434666

$ws.Range("F3").Value = 0.0002084971397132102
$ws.Range("G3").Value = [double]"7.575429228818257e-08"
$ws.Range("H3").Value = 0.0005872800645465599
$ws.Range("I3").Value = 0.0001957429104065394
$ws.Range("J3").Value = [double]"7.245286614021769e-08"
$ws.Range("K3").Value = 0.000551300869807884
$ws.Range("L3").Value = 0.0002091858909706858
$ws.Range("M3").Value = [double]"7.605071481457858e-08"
$ws.Range("N3").Value = 0.0005894172111076503

$ws.Range("F4").Value = 34.32240662895164
$ws.Range("G4").Value = 33.93347633059082
$ws.Range("H4").Value = 34.70549484905773
$ws.Range("I4").Value = 3.833293123652094
$ws.Range("J4").Value = 3.80546947756305
$ws.Range("K4").Value = 3.861071148575839
$ws.Range("L4").Value = 0.2690451030365905
$ws.Range("M4").Value = 0.2667846464728477
$ws.Range("N4").Value = 0.2714705563545743
